# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45192 (2023-09-23) to serial date 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 367; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value2 = 45202
    }
}
